# Target sheet: "ランサーズ" (sheet1) — the scraped Lancers job listing.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Drop every existing hyperlink first so the stale F3/F4/F5 links (attached
# to rows we are about to remove) don't linger in the saved package.
$ws.Cells.Hyperlinks.Delete()

# The old rows 3-5 disappear entirely; only the header + one data row remain.
$ws.Rows("3:5").Delete()

# Row 2 becomes the freshly scraped listing.
$ws.Range("A2").Value = "2026-01-02 01:26:58"
$ws.Range("B2").Value = "複数WEBサイトへの日記一括投稿ツールの修正 or 新規作成をお願いしたいです"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5463948"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5463948")
$ws.Range("G2").Value = 90
$ws.Range("H2").Value = "◆ツール ◇サイト"

# Column widths narrow slightly (B: 51->41, D: 32->22 in saved character
# units). Excel's ColumnWidth property reports/accepts a value that is the
# stored <col width> minus 5/6, so back that out to land on exact integers.
$ws.Columns("B").ColumnWidth = 41 - 5/6
$ws.Columns("D").ColumnWidth = 22 - 5/6
